# Add season record columns (Wins, Losses, Ties) to the roster sheet.
# These three new columns (AD, AE, AF) extend the used range from
# A1:AC63 to A1:AF63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, which uses
# the bold/centered header style) onto the three new header cells so they
# look consistent with the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-63): 67 wins, 95 losses,
# 0 ties.
$ws.Range("AD2:AD63").Value = 67
$ws.Range("AE2:AE63").Value = 95
$ws.Range("AF2:AF63").Value = 0
